$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add two new mapping rows (POI code -> Region), mirroring the existing
# "EG-FYM" rows already present in the sheet.
$ws.Range("A14").Value = "KK"
$ws.Range("B14").Value = "EG-FYM"

$ws.Range("A15").Value = "KomK"
$ws.Range("B15").Value = "EG-FYM"

# Match the formatting used on the other short POI-code cells (A7:A10, A13)
# by copying their style down onto the new cells, including a trailing
# formatted-but-empty cell at A16 (as left behind in the source edit).
$ws.Range("A7").Copy()
$ws.Range("A14:A16").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Reflect the final selection/scroll state from the edit.
$ws.Range("C20").Select()
